# Rename the inline logo pictures in the document's headers/footers.
#
#   * Pearson logo inline pictures (in both footers) go from
#     "image1.png" -> "image2.png"
#   * BTEC logo inline picture (in the header) goes from
#     "image2.jpg" -> "image1.jpg"
#
# The pictures themselves (their embedded binary data / relationship
# targets) are untouched - only the shape's display name changes.

$d = $word.ActiveDocument

function Rename-InlineLogo($Range, $AltText, $NewName) {
    $shapes = $Range.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.AlternativeText -eq $AltText) {
            $shape.Name = $NewName
        }
    }
}

$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

# Walk every section's headers and footers (primary / first-page / even-page)
# so the rename applies no matter which header/footer collection holds the
# pictures.
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $header = $section.Headers.Item($h)
        Rename-InlineLogo $header.Range $pearsonAlt "image2.png"
        Rename-InlineLogo $header.Range $btecAlt "image1.jpg"
    }

    for ($f = 1; $f -le 3; $f++) {
        $footer = $section.Footers.Item($f)
        Rename-InlineLogo $footer.Range $pearsonAlt "image2.png"
        Rename-InlineLogo $footer.Range $btecAlt "image1.jpg"
    }
}
